# Delete the "찬양해 / Sing, sing, sing" slide (original slide 2).
$p = $ppt.ActivePresentation
$p.Slides.Item(2).Delete()

# The slide that used to be slide 3 is now slide 2. Update its body
# placeholder text to add the "Sing, sing, sing" refrain as a new first
# paragraph, and resize/reposition that shape.
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Sing, sing, sing`rAnd make music with the heavens`rWe will sing, sing, sing"
$body.Left = -992 / 914400 * 72
$body.Top = 666658 / 914400 * 72
$body.Width = 12193984 / 914400 * 72
$body.Height = 2762341 / 914400 * 72

# Delete the trailing blank slide (original slide 12).
$p.Slides.Item($p.Slides.Count).Delete()
